$d = $word.ActiveDocument

# ===========================================================================
# Part 1 - Merge runs that exist only because spell-check wrapped a word in
# <w:proofErr>...</w:proofErr>. Re-running Find & Replace across the whole
# phrase (a phrase that already spans the proofErr-wrapped sub-runs) makes
# Word collapse the matched text back into a single run and drop the
# now-unnecessary w:proofErr start/end markers - exactly what the target
# XML looks like.
# ===========================================================================

$d.Content.Find.Execute(
    "There are 8 Teachers: Ava, Isabella, Harper, Elijah, DrSmith, DrJohnson, and 2 instances with missing identifiers (one teaches Chemistry101 and the other Physics101).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "There are 8 Teachers: Ava, Isabella, Harper, Elijah, DrSmith, DrJohnson, and 2 instances with missing identifiers (one teaches Chemistry101 and the other Physics101).",
    2) | Out-Null

$d.Content.Find.Execute(
    "There are 2 Promotions: AcademicExcellenceScholarship and MeritScholarship.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "There are 2 Promotions: AcademicExcellenceScholarship and MeritScholarship.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Main classes include Person, Teacher, Student, EducationalEntity, Program, Course, and Promotion.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Main classes include Person, Teacher, Student, EducationalEntity, Program, Course, and Promotion.",
    2) | Out-Null

# ===========================================================================
# Part 2 - Append seven new commentary paragraphs (one of which carries two
# runs split by a rendered page break) right after the "Overall, the RDF
# data..." paragraph and before the document's final (empty) paragraph.
#
# InsertXML() replaces whatever paragraph the target Range touches, so each
# new paragraph is built in two safe steps:
#   1. InsertParagraphAfter() on the anchor paragraph - this only ever adds
#      a brand-new *empty* paragraph right after the anchor, leaving the
#      anchor paragraph itself untouched.
#   2. InsertXML() against that brand-new empty paragraph's own Range (which
#      does not overlap any neighboring paragraph), replacing its (empty)
#      content with the real paragraph markup.
# The newly written paragraph then becomes the anchor for the next one, so
# the whole block is rebuilt in original order.
# ===========================================================================

$paragraphXmls = @(
@'
<w:p><w:pPr><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>--------------------</w:t></w:r></w:p>
'@,
@'
<w:p><w:pPr><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>This is an RDF graph that represents a university system. It contains information about people (students and teachers), courses, programs, and promotions/scholarships. The graph uses several prefixes to define namespaces for common RDF vocabularies, such as rdf, rdfs, and owl.</w:t></w:r></w:p>
'@,
@'
<w:p><w:pPr><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>The graph defines four main classes: Person, Teacher, Student, EducationalEntity, Program, Course, and Promotion. Person is the superclass of Teacher and Student, and EducationalEntity is the superclass of Program, Course, and Promotion.</w:t></w:r></w:p>
'@,
@'
<w:p><w:pPr><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'>There are several object properties defined in the graph, including teaches, taughtBy, enrolledIn, partOfProgram, and awardedTo. There are also several datatype properties defined in the graph, </w:t></w:r><w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:lastRenderedPageBreak/><w:t>such as emailAddress, dateOfBirth, hasEnrollmentStatus, courseTitle, courseDescription, teacherName, and studentName.</w:t></w:r></w:p>
'@,
@'
<w:p><w:pPr><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>The graph contains information about several individuals, including Ava, Ethan, Liam, Isabella, Amelia, Harper, Daniel, Elijah, Emily, Samuel, Abigail, Alexander, William, Olivia, DrSmith, and DrJohnson. These individuals are instances of either Teacher or Student.</w:t></w:r></w:p>
'@,
@'
<w:p><w:pPr><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>The graph also contains information about several courses, including Math101, Psychology101, Physics101, Chemistry101, Literature101, and History101. Each course has a course title, a course description, and a program it is part of.</w:t></w:r></w:p>
'@,
@'
<w:p><w:pPr><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>The graph defines several programs, including Bachelor of Science in Computer Science, Bachelor of Arts in Psychology, Bachelor of Science in Physics, Bachelor of Science in Chemistry, Bachelor of Arts in Literature, and Bachelor of Arts in History. Each program has a program name and a program duration.</w:t></w:r></w:p>
'@,
@'
<w:p><w:pPr><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>Finally, the graph defines two promotions/scholarships: AcademicExcellenceScholarship and MeritScholarship. Both promotions are awarded to specific students.</w:t></w:r></w:p>
'@
)

# Locate the "Overall, the RDF data..." paragraph to use as the starting
# anchor (found by content rather than a hard-coded index, in case earlier
# edits shifted paragraph numbering).
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Overall, the RDF data and ontology describe an academic setting*") {
        $anchorIndex = $i
    }
}

foreach ($px in $paragraphXmls) {
    $anchor = $d.Paragraphs.Item($anchorIndex)
    $r = $anchor.Range.Duplicate
    $r.Collapse(0)
    $r.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($anchorIndex + 1)
    $newPara.Range.Duplicate.InsertXML($px)

    $anchorIndex = $anchorIndex + 1
}

Write-Output "edit complete"
